# Apply golden-test style updates described by the commit diff.
$d = $word.ActiveDocument

# 1. New paragraph style "AbstractTitle" (display name "Abstract Title"),
#    based on Normal, followed by Abstract, centered/bold/blue heading.
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# 2. "Abstract" style: reduce space-before from 300 (twips/20 => 15pt) to
#    100 (=> 5pt); space-after stays at 300 (15pt).
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 3. "ImportTok" character style gains green, bold formatting.
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# 4. "BuiltInTok" character style gains green formatting.
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768
